$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.291.57"
$ws.Range("E2").Value = "  +6.75%  "

$ws.Range("D3").Value = "3.117.83"
$ws.Range("E3").Value = "  +4.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.64%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.107.94"
$ws.Range("E8").Value = "  +4.35%  "

$ws.Range("E9").Value = "  +2.39%  "

$ws.Range("E10").Value = "  +10.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.87%  "

$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("E13").Value = "  +7.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.19%  "

$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "3.631.09"
$ws.Range("E16").Value = "  +4.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "63.204.53"
$ws.Range("E18").Value = "  +6.58%  "

$ws.Range("D19").Value = "3.113.70"
$ws.Range("E19").Value = "  +4.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("E23").Value = "  +6.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("E28").Value = "  +6.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.63%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.75%  "

$ws.Range("E32").Value = "  +10.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.98%  "

$ws.Range("D35").Value = "0.0₃0819"
$ws.Range("E35").Value = "  +7.86%  "

$ws.Range("E36").Value = "  +4.27%  "

$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "51.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "429.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.67%  "

$ws.Range("D42").Value = "2.971.90"
$ws.Range("E42").Value = "  +7.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0374"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("E45").Value = "  +9.77%  "

$ws.Range("E46").Value = "  +8.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.58%  "

$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.41%  "
